# "Added three new tests to moran's I doc." -- fill in the placeholder
# "TODO" results now that the tests have actually been run, replacing
# them with their real significance markers ("****" = significant at
# 0.001, "x" = insignificant). Three of the five results are on the
# Moran's I row (A6); the other two are the matching Mantel / Mantel
# Correlogram results for the same "88 soils" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Mantel), col O (88 soils / Latitude (distances)) -> highly significant
$ws.Range("O4").Value = "****"

# Row 5 (Mantel Correlogram), col N (88 soils / PH (distances) shuffled) -> insignificant
$ws.Range("N5").Value = "x"

# Row 6 (Moran's I), col K (Glen Canyon / Years since submerged (distances)) -> highly significant
$ws.Range("K6").Value = "****"
# Row 6 (Moran's I), col L (Glen Canyon / Years since submerged (distances) shuffled) -> insignificant
$ws.Range("L6").Value = "x"
# Row 6 (Moran's I), col O (88 soils / Latitude (distances)) -> highly significant
$ws.Range("O6").Value = "****"

# Leave the selection where the author ended up after editing
$ws.Range("O15").Select() | Out-Null
